$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the extra row (previously "Martin") entirely, shifting nothing below it up
$ws.Rows.Item(17).Delete()

# Restore the example/template names in place of the submitted group's names
$ws.Range("B12").Value = "Alice"
$ws.Range("B13").Value = "Bob"
$ws.Range("B14").Value = "Claire"
$ws.Range("B15").Value = "David"
$ws.Range("B16").Value = "Elaine"

# These rows no longer use the special black-font style; restore default formatting
$ws.Range("B12:B16").ClearFormats()

# Match the reported selection after the edit
$ws.Range("E13:E16").Select()
